$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.754.02"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "2.908.26"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.556"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.79%  "
$ws.Range("D9").Value = "2.917.86"
$ws.Range("E9").Value = "  -1.44%  "
$ws.Range("E10").Value = "  -2.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.02"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.365"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").Value = "3.419.36"
$ws.Range("E13").Value = "  -1.41%  "
$ws.Range("E14").Value = "  +1.43%  "
$ws.Range("D15").Value = "60.717.05"
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "2.912.06"
$ws.Range("E17").Value = "  -1.65%  "
$ws.Range("E18").Value = "  -1.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "362.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  -2.34%  "
$ws.Range("E27").Value = "  -1.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.55%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("E32").Value = "  -1.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.41"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.59"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.00%  "
$ws.Range("E37").Value = "  -3.56%  "
$ws.Range("E38").Value = "  -4.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.97%  "
$ws.Range("E40").Value = "  -2.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.41%  "
$ws.Range("D42").Value = "2.288.05"
$ws.Range("E42").Value = "  -4.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.650"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0584"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.58%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.29%  "
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.997"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("E48").Value = "  -2.00%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0927"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.62%  "
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "252.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.37%  "
